$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# Columns A (date-like) and C (numeric-like) must stay plain text, not get
# auto-converted to a date serial / number by Excel's input parsing.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-10-12"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "Pick 3"

$cellC = $ws.Cells.Item($row, 3)
$cellC.NumberFormat = "@"
$cellC.Value = "251012"
$cellC.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "7-0-2"
$ws.Cells.Item($row, 5).Value = "2025-10-12T21:35:21.827+04:00"
